$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Capture Sorted")

$ws.Range("B4").Value = 9
$ws.Range("C4").Value = " VDD_48"
$ws.Range("E4").Value = 16
$ws.Range("K4").Value = "z"
$ws.Range("L4").Value = " CKPWRGD/PD#"
$ws.Range("N4").Value = 16

$ws.Range("B5").Value = 47
$ws.Range("C5").Value = " VDD_CPU"
$ws.Range("E5").Value = 18
$ws.Range("K5").Value = 29
$ws.Range("L5").Value = " CPU_STOP#/SRC_5#"
$ws.Range("N5").Value = 10

$ws.Range("B6").Value = 41
$ws.Range("C6").Value = " VDD_CPU_I/O"
$ws.Range("E6").Value = 14
$ws.Range("K6").Value = 49
$ws.Range("L6").Value = " FS_B/TEST_MODE"
$ws.Range("N6").Value = 17

$ws.Range("B7").Value = 12
$ws.Range("C7").Value = " VDD_I/O"
$ws.Range("E7").Value = 9
$ws.Range("K7").Value = 40
$ws.Range("L7").Value = " I/O_Vout O"
$ws.Range("N7").Value = 2

$ws.Range("B8").Value = 2
$ws.Range("C8").Value = " VDD_PCI"
$ws.Range("E8").Value = 15
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = " PCI_0/CLKREQ_A#"
$ws.Range("N8").Value = 3

$ws.Range("B9").Value = 16
$ws.Range("C9").Value = " VDD_PLL3"
$ws.Range("E9").Value = 17
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = " PCI_1/CLKREQ_B#"
$ws.Range("N9").Value = 4

$ws.Range("B10").Value = 20
$ws.Range("C10").Value = " VDD_PLL3_I/O"
$ws.Range("E10").Value = 10
$ws.Range("K10").Value = 6
$ws.Range("L10").Value = " PCI_4/SRC_5_EN"
$ws.Range("N10").Value = 6

$ws.Range("B11").Value = 53
$ws.Range("C11").Value = " VDD_REF"
$ws.Range("E11").Value = 19
$ws.Range("K11").Value = 30
$ws.Range("L11").Value = " PCI_STOP#/SRC_5"
$ws.Range("N11").Value = 11

$ws.Range("B12").Value = 31
$ws.Range("C12").Value = " VDD_SRC"
$ws.Range("E12").Value = 12
$ws.Range("K12").Value = 7
$ws.Range("L12").Value = " PCIF_5/ITP_EN"
$ws.Range("N12").Value = 7

$ws.Range("B13").Value = 26
$ws.Range("C13").Value = " VDD_SRC_I/O_1"
$ws.Range("E13").Value = 11
$ws.Range("K13").Value = 54
$ws.Range("L13").Value = " REF/FS_C/TEST_SE L"
$ws.Range("N13").Value = 14

$ws.Range("B14").Value = 37
$ws.Range("C14").Value = " VDD_SRC_I/O_2"
$ws.Range("E14").Value = 13
$ws.Range("K14").Value = 56
$ws.Range("L14").Value = " SCL"
$ws.Range("N14").Value = 19

$ws.Range("B15").Value = 11
$ws.Range("C15").Value = " VSS_48"
$ws.Range("E15").Value = 3
$ws.Range("K15").Value = 55
$ws.Range("L15").Value = " SDA"
$ws.Range("N15").Value = 15

$ws.Range("B16").Value = 44
$ws.Range("C16").Value = " VSS_CPU"
$ws.Range("E16").Value = 7
$ws.Range("K16").Value = 25
$ws.Range("L16").Value = " SRC_3#/CLKREQ_D#"
$ws.Range("N16").Value = 9

$ws.Range("B17").Value = 15
$ws.Range("C17").Value = " VSS_I/O"
$ws.Range("E17").Value = 1
$ws.Range("K17").Value = 24
$ws.Range("L17").Value = " SRC_3/CLKREQ_C#"
$ws.Range("N17").Value = 8

$ws.Range("B18").Value = 8
$ws.Range("C18").Value = " VSS_PCI"
$ws.Range("E18").Value = 2
$ws.Range("K18").Value = 35
$ws.Range("L18").Value = " SRC_7#/CLKREQ_E#"
$ws.Range("N18").Value = 12

$ws.Range("B19").Value = 19
$ws.Range("C19").Value = " VSS_PLL3"
$ws.Range("E19").Value = 4
$ws.Range("K19").Value = 36
$ws.Range("L19").Value = " SRC_7/CLKREQ_F#"
$ws.Range("N19").Value = 13

$ws.Range("B20").Value = 50
$ws.Range("C20").Value = " VSS_REF"
$ws.Range("E20").Value = 8
$ws.Range("K20").Value = 4
$ws.Range("L20").Value = " TME/PCI_2"
$ws.Range("N20").Value = 5

$ws.Range("B21").Value = 23
$ws.Range("C21").Value = " VSS_SRC_1"
$ws.Range("E21").Value = 5
$ws.Range("K21").Value = 10
$ws.Range("L21").Value = " USB/FS_A I/O"
$ws.Range("N21").Value = 1

$ws.Range("B22").Value = 34
$ws.Range("C22").Value = " VSS_SRC_2"
$ws.Range("E22").Value = 6
$ws.Range("K22").Value = 52
$ws.Range("L22").Value = " XTAL_IN"
$ws.Range("N22").Value = 18

$ws.Range("K3:O22").AutoFilter()
